# Weekly price-sheet update: insert a new week's record as row 10,
# pushing the existing historical rows (old 10-92) down to (11-93).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10 (shifts rows 10:92 -> 11:93,
# and extends the used range / dimension to A1:R93 automatically).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with this week's data.
$ws.Cells.Item(10, 1).Value  = 9
$ws.Cells.Item(10, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(10, 3).Value  = "Metropolitana"
$ws.Cells.Item(10, 4).Value  = 44545
$ws.Cells.Item(10, 5).Value  = 13
$ws.Cells.Item(10, 6).Value  = 100112022
$ws.Cells.Item(10, 7).Value  = "Arveja Verde"
$ws.Cells.Item(10, 8).Value  = "Sin especificar"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 34
$ws.Cells.Item(10, 11).Value = 17000
$ws.Cells.Item(10, 12).Value = 18000
$ws.Cells.Item(10, 13).Value = 17500
$ws.Cells.Item(10, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(10, 15).Value = "Carahue"
$ws.Cells.Item(10, 16).Value = 700
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
